$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1234
$ws1.Range("F4").Value = 16988
$ws1.Range("F5").Value = 41
$ws1.Range("F13").Value = 11760
$ws1.Range("F15").Value = 21
$ws1.Range("F16").Value = 1445
$ws1.Range("F17").Value = 4675
$ws1.Range("F18").Value = 483
$ws1.Range("F20").Value = 410
$ws1.Range("F21").Value = 75

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1234
$ws4.Range("F5").Value = 16988
$ws4.Range("F6").Value = 41
$ws4.Range("F14").Value = 4
$ws4.Range("F16").Value = 11760
$ws4.Range("F18").Value = 21
$ws4.Range("F19").Value = 1445
$ws4.Range("F20").Value = 4675
$ws4.Range("F21").Value = 483
$ws4.Range("F23").Value = 410
$ws4.Range("F24").Value = 75
